# Splits the "E-mail: cursods_profdanilo@gmail.com" run into three runs:
#   "E-mail: " + "cursods.profdanilo" + "@gmail.com"
# (the middle segment had its underscore turned into a dot).

$p = $ppt.ActivePresentation

$oldMiddle = "cursods_profdanilo"
$newMiddle = "cursods.profdanilo"
$fullOld   = "E-mail: " + $oldMiddle + "@gmail.com"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)

            if ($para.Text.Contains($fullOld)) {
                $startOfMiddle = $para.Text.IndexOf($oldMiddle) + 1
                $midRange = $para.Characters($startOfMiddle, $oldMiddle.Length)
                $midRange.Text = $newMiddle
            }
        }
    }
}
